$wb = $excel.ActiveWorkbook

# Parametric survival model output sheets: update est (B) / se (C) for
# each distribution's parameters (row 2 / row 3).
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -3.09126802584416
$ws.Range("C2").Value = 0.175312283413101
$ws.Range("B3").Value = 0.143473106550909
$ws.Range("C3").Value = 0.117818399725227

$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 2.46833408038063
$ws.Range("C2").Value = 0.22714291112948
$ws.Range("B3").Value = -1.03846010749416
$ws.Range("C3").Value = 0.115243977650148

$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -2.3370003112681
$ws.Range("C2").Value = 0.129544185691206
$ws.Range("B3").Value = 1.73953940198937
$ws.Range("C3").Value = 0.208695357802345

$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.66618220861419
$ws.Range("C2").Value = 0.144808633201942
$ws.Range("B3").Value = -0.0198690908444268
$ws.Range("C3").Value = 0.0191775399753587

# Covariance matrices (2x2) for each distribution's parameters.
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0307343967155154
$ws.Range("B2").Value = -0.0141763938303875
$ws.Range("A3").Value = -0.0141763938303875
$ws.Range("B3").Value = 0.0138811753138134

$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.0515939020763749
$ws.Range("B2").Value = -0.0232637853881486
$ws.Range("A3").Value = -0.0232637853881486
$ws.Range("B3").Value = 0.0132811743846279

$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.0167816960463977
$ws.Range("B2").Value = 0.0130353650991837
$ws.Range("A3").Value = 0.0130353650991837
$ws.Range("B3").Value = 0.0435537523682489

$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0209695402498145
$ws.Range("B2").Value = -0.00167073607669963
$ws.Range("A3").Value = -0.00167073607669963
$ws.Range("B3").Value = 0.00036777803950648
